$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-order the header columns in row 2.
# Old order: BusinessKey, Code, Framework_ID, Name, OrganizationBusinessKey
# New order: Framework_ID, BusinessKey, OrganizationBusinessKey, Code, Name
$ws.Range("A2").Value = "Framework_ID"
$ws.Range("B2").Value = "BusinessKey"
$ws.Range("C2").Value = "OrganizationBusinessKey"
$ws.Range("D2").Value = "Code"
$ws.Range("E2").Value = "Name"
